# Rename the header row suffixes from the generic "_old" / "_new" markers
# to the concrete input-file-version markers "_FV2310" / "_FV2404", then
# turn the header row into a proper Excel Table and freeze it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A-J ("…_old") -> "…_FV2310"
$fv2310Headers = @(
  "Segmentname_FV2310",
  "Segmentgruppe_FV2310",
  "Segment_FV2310",
  "Datenelement_FV2310",
  "Segment ID_FV2310",
  "Code_FV2310",
  "Qualifier_FV2310",
  "Beschreibung_FV2310",
  "Bedingungsausdruck_FV2310",
  "Bedingung_FV2310"
)

# Column K ("diff") stays untouched.

# Columns L-U ("…_new") -> "…_FV2404"
$fv2404Headers = @(
  "Segmentname_FV2404",
  "Segmentgruppe_FV2404",
  "Segment_FV2404",
  "Datenelement_FV2404",
  "Segment ID_FV2404",
  "Code_FV2404",
  "Qualifier_FV2404",
  "Beschreibung_FV2404",
  "Bedingungsausdruck_FV2404",
  "Bedingung_FV2404"
)

for ($i = 0; $i -lt $fv2310Headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $fv2310Headers[$i]
}
for ($i = 0; $i -lt $fv2404Headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $fv2404Headers[$i]
}

# Turn the whole used range into a native Excel Table ("Table1") so headers
# get filter buttons / structured references, matching the new xl/tables/table1.xml part.
$tableRange = $ws.Range("A1:U74")
$table = $ws.ListObjects.Add(1, $tableRange, [System.Reflection.Missing]::Value, 1)
$table.Name = "Table1"
$table.TableStyle = ""

# Freeze the header row (row 1) so it stays visible while scrolling.
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
